# Add network configuration for Project List file access and update invoice pipeline
# Inserts a new invoice row (Norman S. Wright Mech. Equip. LLC.) at the top of the
# data block, pushing the existing rows down by one, refreshes the per-row
# "file_name" re-export stamps for the shifted rows, and widens the
# Vendor_Name column to fit the new, longer vendor name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a brand-new row right under the header row and strip any
#     formatting it might have inherited from the row above it. ---
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# --- Populate the newly inserted row 2 with the Norman S. Wright invoice. ---
$ws.Range("A2").Value = "AA1"
$ws.Range("B2").Value = "08/13/25"
$ws.Range("C2").Value = "2025-08-22"
$ws.Range("D2").Value = "NORWRI"
$ws.Range("E2").Value = "Norman S. Wright Mech. Equip. LLC."
$ws.Range("F2").Value = "I"
$ws.Range("G2").Value = "127813"
$ws.Range("H2").Value = "08/13/25"
$ws.Range("I2").Value = "821.02"
$ws.Range("J2").Value = "73.79"
$ws.Range("L2").Value = "747.23"
$ws.Range("M2").Value = 1466
$ws.Range("P2").Value = "5040"
$ws.Range("Q2").Value = 330
$ws.Range("R2").Value = "E"
$ws.Range("T2").Value = "127813_1755894601049.pdf"

# --- The rows that were pushed down keep their original data, but each was
#     re-exported with a new processing timestamp, so refresh file_name. ---
$ws.Range("T3").Value = "1_1755894595697.pdf"
$ws.Range("T4").Value = "2_1755894595700.pdf"
$ws.Range("T5").Value = "3_1755894595693.pdf"
$ws.Range("T6").Value = "4_1755894595691.pdf"
$ws.Range("T7").Value = "inv-01-875854.pdf_page_1_1755894595704.pdf"
$ws.Range("T8").Value = "sin221250_page_1_1755894595706.pdf"

# --- Widen the Vendor_Name column so the longer vendor name isn't clipped. ---
$ws.Range("E1").ColumnWidth = 35.14
